$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iEX_TC_ID_1"
$ws.Range("B2").Value = "@Smoke Verify Elumina Login and Create Exam"
$ws.Range("C2").Value = "interrupted"
